# Minor fix to presentation
# Adds a new "Objects" bullet (level 2 / lvl="1", same 20pt size as its
# siblings) right before the existing "User Interface" bullet on the
# "Business & System Objectives" slide's content placeholder.

$p = $ppt.ActivePresentation

# Locate the slide/shape that contains the "User Interface" bullet instead
# of hard-coding a slide index, so the script is resilient to reordering.
$targetShape = $null

foreach ($s in $p.Slides) {
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -like "*User Interface*") {
                $targetShape = $shp
            }
        }
    }
}

if ($targetShape -eq $null) {
    throw "Could not find the 'User Interface' bullet in the presentation."
}

$tr = $targetShape.TextFrame.TextRange

# Find the "User Interface" paragraph so the new bullet can be inserted
# immediately before it (i.e. right after "Build System"). Paragraph text
# includes a trailing carriage return, so trim it before comparing.
$paraCount = $tr.Paragraphs().Count
$uiIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $paraText = $para.Text.TrimEnd([char]13)
    if ($paraText -eq "User Interface") {
        $uiIndex = $i
    }
}

if ($uiIndex -eq -1) {
    throw "Could not find the 'User Interface' paragraph."
}

$uiPara = $tr.Paragraphs($uiIndex, 1)

# Insert a new paragraph ("Objects") before "User Interface"; it naturally
# inherits the same level (lvl="1") and font size (20pt) as the paragraph
# it is being inserted in front of.
$newRange = $uiPara.InsertBefore("Objects" + [char]13)
